# Assignment1.xlsx update: "Add files via upload"
# Applies numeric cell updates on the "PO List" sheet and updates the
# active selection on that sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PO List")

# --- Row 3 ---
$ws.Range("N3").Value = 11
$ws.Range("R3").Value = 12

# --- Row 4 ---
$ws.Range("N4").Value = 18
$ws.Range("R4").Value = 6

# --- Row 5 ---
$ws.Range("N5").Value = 22
$ws.Range("R5").Value = 2

# --- Row 6 ---
$ws.Range("N6").Value = 15
$ws.Range("R6").Value = 17

# --- Row 7 ---
$ws.Range("N7").Value = 24
$ws.Range("R7").Value = 16

# --- Row 8 ---
$ws.Range("N8").Value = 21

# --- Row 9 ---
$ws.Range("N9").Value = 4
$ws.Range("R9").Value = 14

# --- Row 10 ---
$ws.Range("N10").Value = 10
$ws.Range("R10").Value = 11

# --- Row 11 ---
$ws.Range("N11").Value = 8
$ws.Range("R11").Value = 20

# --- Row 12 ---
$ws.Range("N12").Value = 6
$ws.Range("R12").Value = 18

# --- Row 13 ---
$ws.Range("N13").Value = 5

# --- Row 14 ---
$ws.Range("N14").Value = 20

# --- Row 15 ---
$ws.Range("N15").Value = 12
$ws.Range("R15").Value = 8

# --- Row 16 ---
$ws.Range("N16").Value = 14
$ws.Range("R16").Value = 9

# --- Row 17 ---
$ws.Range("N17").Value = 13
$ws.Range("R17").Value = 21

# --- Row 18 ---
$ws.Range("N18").Value = 23
$ws.Range("R18").Value = 4

# --- Row 19 ---
$ws.Range("N19").Value = 2

# --- Row 20 ---
$ws.Range("N20").Value = 19
$ws.Range("R20").Value = 10

# --- Row 21 ---
$ws.Range("N21").Value = 17
$ws.Range("R21").Value = 15

# --- Row 23 ---
$ws.Range("R23").Value = 19

# --- Row 24 ---
$ws.Range("N24").Value = 7
$ws.Range("R24").Value = 13

# --- Row 25 ---
$ws.Range("N25").Value = 3
$ws.Range("O25").Value = 3
$ws.Range("P25").Value = 3
$ws.Range("Q25").Value = 44918
$ws.Range("R25").Value = 1

# --- Row 26 ---
$ws.Range("N26").Value = 15
$ws.Range("R26").Value = 7

# --- Row 27 ---
$ws.Range("K27").Value = 3
$ws.Range("L27").Value = 3
$ws.Range("M27").Value = 44914
$ws.Range("N27").Value = 1
$ws.Range("R27").Value = 5

# --- Row 29 ---
$ws.Range("N29").Value = 9
$ws.Range("R29").Value = 2

# --- Update the active selection on the sheet (bottomRight pane) ---
$ws.Activate()
$ws.Range("P40").Select()
